# feat: add 2022-Q3 data
#
# Inserts a new worksheet "2022-Q3" right after the "总计" (summary) sheet and
# before the existing "2022-Q2" sheet, populates it with the new quarter's
# fund-holding data, and adds a corresponding summary row at the top of the
# "总计" sheet's data table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q3" worksheet immediately before "2022-Q2".
# ---------------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($refSheet)
$newSheet.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# 2) Populate headers (row 1, columns B:H) with the same bold/boxed style
#    used by every other quarterly sheet.
# ---------------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$headerCols = @("B", "C", "D", "E", "F", "G", "H")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $newSheet.Range($headerCols[$i] + "1")
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# ---------------------------------------------------------------------------
# 3) Populate the fund rows. Column A is the numeric pandas-style row index
#    (bold/boxed like the header), column H is a plain number (rank); all of
#    B/D/E/F/G must stay TEXT (fund codes keep leading zeros, ratios keep
#    their exact printed form) even though they look numeric, so they are
#    entered with a leading apostrophe to force text storage. C is already
#    non-numeric text and needs no special handling.
# ---------------------------------------------------------------------------
$rows = @(
    @("003713", "英大睿盛灵活配置混合A", "2.83", "93.65", "5.78", "0.1636", 10),
    @("003714", "英大睿盛灵活配置混合C", "2.19", "93.65", "5.78", "0.1266", 10),
    @("014179", "中银证券远见价值混合A", "1.56", "93.65", "3.77", "0.0588", 7),
    @("001607", "英大策略优选混合A", "0.57", "91.98", "4.58", "0.0261", 9),
    @("014180", "中银证券远见价值混合C", "0.16", "93.65", "3.77", "0.0060", 7),
    @("562530", "华夏中证智选1000价值稳健策略ETF", "0.54", "94.32", "0.99", "0.0053", 1),
    @("001608", "英大策略优选混合C", "0.02", "91.98", "4.58", "0.0009", 9)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $idxCell = $newSheet.Range("A" + $r)
    $idxCell.Value = $i
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    $newSheet.Range("B" + $r).Value = "'" + $row[0]
    $newSheet.Range("C" + $r).Value = $row[1]
    $newSheet.Range("D" + $r).Value = "'" + $row[2]
    $newSheet.Range("E" + $r).Value = "'" + $row[3]
    $newSheet.Range("F" + $r).Value = "'" + $row[4]
    $newSheet.Range("G" + $r).Value = "'" + $row[5]
    $newSheet.Range("H" + $r).Value = $row[6]
}

# ---------------------------------------------------------------------------
# 4) Insert a new row 2 at the top of the "总计" sheet's data table for the
#    new "2022-Q3" quarter, shifting the existing B:D rows down by one.
#    Column A is a plain 0-based row counter (0,1,2,...) independent of
#    which quarter occupies the row, so after the shift it is rewritten as
#    a fixed 0..6 sequence rather than being carried along with the insert.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# Clear the blank formatting the row-insert left on B2:D2 so the new data
# row matches the plain (unstyled) look of every other data row.
$summary.Range("B2:D2").ClearFormats()

# Copy column A's box/bold style from the row below (now row 3, an existing
# data row) onto the new A2 index cell.
$summary.Range("A3").Copy() | Out-Null
$summary.Range("A2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 7
$summary.Range("D2").Value = 0.39

# Re-stamp the 0-based row counter in column A for every data row (2-8):
# it always runs 0..6 top-to-bottom regardless of the quarter shuffle above.
for ($r = 2; $r -le 8; $r++) {
    $summary.Range("A" + $r).Value = $r - 2
}
